# Append a new transactional-data row (Order ID / Order Date / Order Status)
# to the bottom of the Orders sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order ID and Order Date must be stored as plain text (like every other
# row above), not auto-coerced into a number/date serial. Leading the
# literal with an apostrophe forces Excel to store it as text; resetting
# the cell style back to "Normal" afterwards drops the quote-prefix
# formatting so the cell ends up with the same default style as its
# neighbours.
$ws.Range("A11").Value = "'1958442"
$ws.Range("A11").Style = "Normal"

$ws.Range("B11").Value = "'04/16/2025"
$ws.Range("B11").Style = "Normal"

$ws.Range("C11").Value = "Pending"
